# Refresh the crypto price table to the latest scrape.
# Every data cell in this sheet is stored as text (Coin/Link/Price/Volume
# columns are all strings, including Price values like "353.40"), so
# plain numeric-looking Price updates are written with the cell
# pre-formatted as Text ("@") to stop Excel from silently coercing them
# into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin / Link / Volume(1h) text updates, plus Price values that are not
# bare numbers (so Excel leaves them as text on their own).
$plainUpdates = @(
    @("D2", "52.142.98"),
    @("D3", "2.942.51"),
    @("E3", "  +4.52%  "),
    @("E4", "  +0.15%  "),
    @("E5", "  +0.61%  "),
    @("E6", "  -0.96%  "),
    @("E7", "  -0.22%  "),
    @("E8", "  -0.02%  "),
    @("E9", "  +0.88%  "),
    @("E10", "  -1.68%  "),
    @("E11", "  +3.38%  "),
    @("E12", "  +1.01%  "),
    @("E13", "  -0.35%  "),
    @("B14", "Polkadot"),
    @("C14", "https://coinranking.com/coin/25W7FG7om+polkadot-dot"),
    @("E14", "  -0.29%  "),
    @("B15", "WrappedliquidstakedEther2.0"),
    @("C15", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"),
    @("D15", "3.404.23"),
    @("E15", "  +4.65%  "),
    @("D16", "2.941.89"),
    @("E16", "  +4.26%  "),
    @("E17", "  +1.07%  "),
    @("D18", "52.234.57"),
    @("E18", "  +0.77%  "),
    @("E19", "  +0.73%  "),
    @("E20", "  -2.85%  "),
    @("E21", "  +5.66%  "),
    @("D22", "0.0₃0981"),
    @("E22", "  +0.45%  "),
    @("E23", "  +0.76%  "),
    @("E24", "  -0.02%  "),
    @("E25", "  +1.21%  "),
    @("E26", "  +11.56%  "),
    @("E27", "  +3.01%  "),
    @("E28", "  +0.06%  "),
    @("E29", "  +12.82%  "),
    @("E30", "  +0.94%  "),
    @("E31", "  +15.17%  "),
    @("B32", "Toncoin"),
    @("C32", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"),
    @("E32", "  -0.30%  "),
    @("B33", "InjectiveProtocol"),
    @("C33", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"),
    @("E33", "  -4.72%  "),
    @("B34", "RenderToken"),
    @("C34", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"),
    @("E34", "  +5.44%  "),
    @("B35", "OKB"),
    @("C35", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"),
    @("E35", "  +0.66%  "),
    @("E36", "  +0.19%  "),
    @("E37", "  -0.18%  "),
    @("E38", "  +5.85%  "),
    @("E40", "  +2.40%  "),
    @("E41", "  +4.42%  "),
    @("E42", "  +1.33%  "),
    @("E43", "  +6.18%  "),
    @("E44", "  -1.59%  "),
    @("D45", "2.199.15"),
    @("E45", "  +2.36%  "),
    @("B46", "NEARProtocol"),
    @("C46", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"),
    @("E46", "  -0.53%  "),
    @("B47", "ApeXProtocol"),
    @("C47", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"),
    @("E47", "  +0.96%  "),
    @("E48", "  -6.10%  "),
    @("E49", "  +10.07%  "),
    @("E50", "  +8.80%  "),
    @("E51", "  -3.66%  ")
)

foreach ($pair in $plainUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Price values that look like bare numbers (e.g. "353.40") - force Text
# format first so Excel keeps them as strings instead of converting to
# numeric cells.
$textPriceUpdates = @(
    @("D5", "353.40"),
    @("D6", "112.18"),
    @("D7", "0.560"),
    @("D9", "0.627"),
    @("D10", "39.57"),
    @("D11", "0.0879"),
    @("D12", "0.137"),
    @("D13", "20.09"),
    @("D14", "7.77"),
    @("D17", "0.984"),
    @("D19", "7.67"),
    @("D20", "3.30"),
    @("D21", "14.25"),
    @("D23", "71.19"),
    @("D24", "268.65"),
    @("D25", "2.80"),
    @("D26", "0.181"),
    @("D27", "27.05"),
    @("D29", "7.03"),
    @("D30", "10.63"),
    @("D31", "0.104"),
    @("D32", "2.26"),
    @("D33", "37.05"),
    @("D34", "6.08"),
    @("D35", "53.09"),
    @("D36", "0.0454"),
    @("D39", "18.70"),
    @("D41", "2.69"),
    @("D43", "23.45"),
    @("D46", "3.50"),
    @("D47", "2.50"),
    @("D48", "114.22"),
    @("D49", "0.249"),
    @("D50", "0.0350"),
    @("D51", "0.953")
)

foreach ($pair in $textPriceUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
}

